{"js": "// Moving from 2.0.2 to 2.0.3.\n// The REF-field runs that used to carry rsidR=\"28A6C2EA6BE54955BE1BCB9DF51B0896\"\n// now carry rsidR=\"42D0D12D881F49E5911F278501AC37AA\", and the bookmark \"Art1\"\n// (bookmarkStart/bookmarkEnd) gets a new w:id.\n\nconst NEW_RSID = \"42D0D12D881F49E5911F278501AC37AA\";\nconst NEW_BOOKMARK_ID = \"160706569393042178039163456731701538084\";\n\nconst NS =\n  'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapPackage(paraXml) {\n  return (\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document ' + NS + \"><w:body>\" +\n    paraXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Rebuild a field-code paragraph (the \"REF Art1 \\h\" field displaying\n// \"Artifact1\"), swapping the rsidR GUID used on its five runs, while\n// preserving the paragraph's own pPr-relevant rsid* attributes.\nfunction fieldParagraphXml(pPrAttrs) {\n  return (\n    \"<w:p\" + pPrAttrs + \">\" +\n    \"<w:r><w:rPr><w:color w:themeColor=\\\"accent6\\\" w:themeShade=\\\"BF\\\" w:val=\\\"E36C0A\\\"/></w:rPr><w:t/></w:r>\" +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF Art1 \\\\h </w:instrText></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>Artifact1</w:t></w:r>' +\n    '<w:r w:rsidR=\"' + NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n    \"</w:p>\"\n  );\n}\n\n// Rebuild the \"Definition of Artifact1\" paragraph, swapping the\n// bookmarkStart/bookmarkEnd w:id for the Art1 bookmark.\nfunction bookmarkParagraphXml(pPrAttrs) {\n  return (\n    \"<w:p\" + pPrAttrs + \">\" +\n    '<w:r w:rsidR=\"00E61FB8\"><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t/></w:r>' +\n    '<w:bookmarkStart w:name=\"Art1\" w:id=\"' + NEW_BOOKMARK_ID + '\"/>' +\n    '<w:r w:rsidR=\"00E61FB8\"><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>Definition of Artifact1</w:t></w:r>' +\n    '<w:bookmarkEnd w:id=\"' + NEW_BOOKMARK_ID + '\"/>' +\n    \"</w:p>\"\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\n// The known (fixed) paragraph-mark attributes of the two field\n// paragraphs and the bookmark paragraph, taken from the document.\nconst FIELD_PPR_ATTRS = [\n  ' w:rsidP=\"00E8765B\" w:rsidR=\"00E8765B\" w:rsidRDefault=\"00E8765B\"',\n  ' w:rsidP=\"00F5495F\" w:rsidR=\"00730F00\" w:rsidRDefault=\"00730F00\"',\n];\nconst BOOKMARK_PPR_ATTRS =\n  ' w:rsidP=\"00F5495F\" w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\"';\n\nlet fieldParaSeen = 0;\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  const range = p.getRange();\n\n  if (text === \"Artifact1\") {\n    const pPrAttrs =\n      FIELD_PPR_ATTRS[fieldParaSeen] || FIELD_PPR_ATTRS[FIELD_PPR_ATTRS.length - 1];\n    fieldParaSeen++;\n    const xml = wrapPackage(fieldParagraphXml(pPrAttrs));\n    range.insertOoxml(xml, Word.InsertLocation.replace);\n  } else if (text === \"Definition of Artifact1\") {\n    const xml = wrapPackage(bookmarkParagraphXml(BOOKMARK_PPR_ATTRS));\n    range.insertOoxml(xml, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Moving from 2.0.2 to 2.0.3.\n# The REF-field runs that used to carry rsidR=\"28A6C2EA6BE54955BE1BCB9DF51B0896\"\n# now carry rsidR=\"42D0D12D881F49E5911F278501AC37AA\", and the bookmark \"Art1\"\n# (bookmarkStart/bookmarkEnd) gets a new w:id.\n\n$NEW_RSID = \"42D0D12D881F49E5911F278501AC37AA\"\n$NEW_BOOKMARK_ID = \"160706569393042178039163456731701538084\"\n\n$NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\nfunction Wrap-Package($paraXml) {\n    return '<?xml version=\"1.0\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document ' + $NS + '><w:body>' +\n        $paraXml +\n        '</w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\nfunction FieldParagraph-Xml($pPrAttrs) {\n    return '<w:p' + $pPrAttrs + '>' +\n        '<w:r><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t/></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:instrText xml:space=\"preserve\"> REF Art1 \\h </w:instrText></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:rPr><w:b w:val=\"true\"/><w:noProof/></w:rPr><w:t>Artifact1</w:t></w:r>' +\n        '<w:r w:rsidR=\"' + $NEW_RSID + '\"><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n        '</w:p>'\n}\n\nfunction BookmarkParagraph-Xml($pPrAttrs) {\n    return '<w:p' + $pPrAttrs + '>' +\n        '<w:r w:rsidR=\"00E61FB8\"><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t/></w:r>' +\n        '<w:bookmarkStart w:name=\"Art1\" w:id=\"' + $NEW_BOOKMARK_ID + '\"/>' +\n        '<w:r w:rsidR=\"00E61FB8\"><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>Definition of Artifact1</w:t></w:r>' +\n        '<w:bookmarkEnd w:id=\"' + $NEW_BOOKMARK_ID + '\"/>' +\n        '</w:p>'\n}\n\n# Known (fixed) paragraph-mark attributes for the two field paragraphs\n# (first and second \"REF Art1\" field) and for the bookmark paragraph.\n$FieldPPrAttrs = @(\n    ' w:rsidP=\"00E8765B\" w:rsidR=\"00E8765B\" w:rsidRDefault=\"00E8765B\"',\n    ' w:rsidP=\"00F5495F\" w:rsidR=\"00730F00\" w:rsidRDefault=\"00730F00\"'\n)\n$BookmarkPPrAttrs = ' w:rsidP=\"00F5495F\" w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\"'\n\n$d = $word.ActiveDocument\n\n$fieldParaSeen = 0\n# Snapshot paragraphs first since replacing content shifts the collection.\n$paraList = @()\nforeach ($p in $d.Paragraphs) {\n    $paraList += $p\n}\n\nforeach ($p in $paraList) {\n    $text = $p.Range.Text\n    if ($text -eq \"Artifact1`r\" -or $text -eq \"Artifact1\") {\n        $idx = [Math]::Min($fieldParaSeen, $FieldPPrAttrs.Length - 1)\n        $pPrAttrs = $FieldPPrAttrs[$idx]\n        $fieldParaSeen = $fieldParaSeen + 1\n        $xml = Wrap-Package (FieldParagraph-Xml $pPrAttrs)\n        [void]$p.Range.InsertXML($xml)\n    } elseif ($text -eq \"Definition of Artifact1`r\" -or $text -eq \"Definition of Artifact1\") {\n        $xml = Wrap-Package (BookmarkParagraph-Xml $BookmarkPPrAttrs)\n        [void]$p.Range.InsertXML($xml)\n    }\n}\n"}
